$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.142.10"
Set-TextCell "E2" "  +1.19%  "
Set-TextCell "D3" "3.906.09"
Set-TextCell "E3" "  -0.98%  "
Set-TextCell "E4" "  -0.10%  "
Set-TextCell "D5" "489.29"
Set-TextCell "E5" "  +3.79%  "
Set-TextCell "D6" "146.41"
Set-TextCell "E6" "  -0.11%  "
Set-TextCell "D7" "0.620"
Set-TextCell "E7" "  -1.25%  "
Set-TextCell "D8" "0.998"
Set-TextCell "E8" "  -0.11%  "
Set-TextCell "D9" "0.727"
Set-TextCell "E9" "  -1.30%  "
Set-TextCell "E10" "  -1.16%  "
Set-TextCell "E11" "  +0.92%  "
Set-TextCell "D12" "42.85"
Set-TextCell "E12" "  -1.67%  "
Set-TextCell "D13" "10.74"
Set-TextCell "E13" "  +2.78%  "
Set-TextCell "D14" "4.520.37"
Set-TextCell "E14" "  -1.05%  "
Set-TextCell "D15" "3.903.01"
Set-TextCell "E15" "  -0.37%  "
Set-TextCell "D16" "14.06"
Set-TextCell "E16" "  -7.17%  "
Set-TextCell "E17" "  -1.18%  "
Set-TextCell "D18" "19.79"
Set-TextCell "E18" "  -0.60%  "
Set-TextCell "E19" "  -2.78%  "
Set-TextCell "D20" "68.216.80"
Set-TextCell "E20" "  +0.98%  "
Set-TextCell "D21" "428.91"
Set-TextCell "E21" "  -1.73%  "
Set-TextCell "D22" "3.55"
Set-TextCell "E22" "  +4.24%  "
Set-TextCell "D23" "14.96"
Set-TextCell "E23" "  +2.95%  "
Set-TextCell "D24" "87.27"
Set-TextCell "E24" "  -0.48%  "
Set-TextCell "D25" "11.33"
Set-TextCell "E25" "  +15.13%  "
Set-TextCell "D26" "11.33"
Set-TextCell "E26" "  +9.77%  "
Set-TextCell "D27" "3.62"
Set-TextCell "E27" "  +0.06%  "
Set-TextCell "D28" "38.14"
Set-TextCell "E28" "  -1.68%  "
Set-TextCell "D29" "5.73"
Set-TextCell "E29" "  -0.36%  "
Set-TextCell "D30" "723.35"
Set-TextCell "E30" "  +0.41%  "
Set-TextCell "D31" "13.77"
Set-TextCell "E31" "  +1.69%  "
Set-TextCell "E32" "  -2.28%  "
Set-TextCell "E33" "  +3.19%  "
Set-TextCell "D34" "6.25"
Set-TextCell "E34" "  +16.53%  "
Set-TextCell "D35" "41.74"
Set-TextCell "E35" "  -2.42%  "
Set-TextCell "D36" "0.0₃0866"
Set-TextCell "E36" "  +7.14%  "
Set-TextCell "D37" "60.47"
Set-TextCell "E37" "  +4.29%  "
Set-TextCell "D38" "0.406"
Set-TextCell "E38" "  +19.18%  "
Set-TextCell "E39" "  -2.10%  "
Set-TextCell "E40" "  -0.10%  "
Set-TextCell "D41" "2.96"
Set-TextCell "E41" "  +14.94%  "
Set-TextCell "D42" "0.0478"
Set-TextCell "E42" "  +0.04%  "
Set-TextCell "D43" "3.14"
Set-TextCell "E43" "  +2.47%  "
Set-TextCell "E44" "  +2.83%  "
Set-TextCell "D45" "0.141"
Set-TextCell "E45" "  -1.28%  "
Set-TextCell "E46" "  +0.00%  "
Set-TextCell "D47" "3.33"
Set-TextCell "E47" "  +4.57%  "
Set-TextCell "D48" "3.40"
Set-TextCell "E48" "  -4.53%  "
Set-TextCell "D49" "2.13"
Set-TextCell "E49" "  -3.62%  "

# Rows 50/51: Monero and BabyDogeCoin swap position with updated values
Set-TextCell "B50" "BabyDogeCoin"
Set-TextCell "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D50" "0.0₆0340"
Set-TextCell "E50" "  +26.45%  "

Set-TextCell "B51" "Monero"
Set-TextCell "C51" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D51" "144.26"
Set-TextCell "E51" "  -2.24%  "
